$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Topic" values (column C) between rows 5-6 and rows 7-8
$c5 = $ws.Range("C5").Value2
$c6 = $ws.Range("C6").Value2
$c7 = $ws.Range("C7").Value2
$c8 = $ws.Range("C8").Value2

$ws.Range("C5").Value = $c7
$ws.Range("C6").Value = $c8
$ws.Range("C7").Value = $c5
$ws.Range("C8").Value = $c6

# Update the selected cell to C16
$ws.Range("C16").Select()
